$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1013, shifting the existing rows 1013:1074 down to 1014:1075
$ws.Rows.Item(1013).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(1013, 1).Value = 5
$ws.Cells.Item(1013, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1013, 3).Value = "Maule"
$ws.Cells.Item(1013, 4).Value = 45267
$ws.Cells.Item(1013, 5).Value = 7
$ws.Cells.Item(1013, 6).Value = "Fruta"
$ws.Cells.Item(1013, 7).Value = 100102
$ws.Cells.Item(1013, 8).Value = "Cítricos"
$ws.Cells.Item(1013, 9).Value = 100102005
$ws.Cells.Item(1013, 10).Value = "Naranja"
$ws.Cells.Item(1013, 11).Value = "Valencia"
$ws.Cells.Item(1013, 12).Value = "Primera"
$ws.Cells.Item(1013, 13).Value = 400
$ws.Cells.Item(1013, 14).Value = 11000
$ws.Cells.Item(1013, 15).Value = 11000
$ws.Cells.Item(1013, 16).Value = 11000
$ws.Cells.Item(1013, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(1013, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1013, 19).Value = 611
$ws.Cells.Item(1013, 20).Value = 18
